$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate 'baselineDate' variable row entirely (row 39), shifting subsequent rows up.
$ws.Rows.Item(39).Delete()

# Add the Coding (column D) for smokingCat (row 29).
$ws.Range("D29").Value = "Current smoker >19 pack year,  Current smoker 1-19 pack year, Current smoker unknown pack year,`nNever smoker,`nPast smoker >19 pack year,`nPast smoker 1-19 pack year,`nPast smoker unknown pack year "

# Add the Coding (column D) for pm10_2010per10 (row 31).
$ws.Range("D31").Value = " "

# Add the Coding (column D) for cancerDate (row 40, after the row shift).
$ws.Range("D40").Value = "   "

# Update the Data Generation column for t_lungCancer (row 48) and cancerDate_Lung (row 49).
$ws.Range("E48").Value = "Caitlyn created column in solid_all_0912.csv. Duplicate eid are filtered out in the AP_data_processing.R script, making this time to first lung cancer diagnosis."
$ws.Range("E49").Value = "Caitlyn created column in solid_all_0912.csv. Duplicate eid are filtered out in the AP_data_processing.R script, making this date of first lung cancer diagnosis."

# Match the saved viewport/selection state (scrolled so row 29 is the top row, with B47 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B47").Select()

Write-Output "done"
